# Update the "date" column (B) from 2021 to 2022 for the "Population"
# indicator rows (878-1142) in the econ_metadata workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B878:B1142").Value = 2022
